# Adds a new "2022-Q4" quarterly sheet right after "总计", and updates the
# "总计" summary sheet with a new row for the 2022-Q4 quarter.
#
# All the other quarterly sheets ("2021-Q3","2021-Q2","2021-Q1","2020-Q4")
# keep their name/content unchanged - they just shift one tab to the right
# to make room for the new sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet immediately after "总计" (i.e.
#    immediately before the current "2021-Q3" sheet), by duplicating the
#    "2021-Q3" sheet so that it starts out with the same layout/styles,
#    then overwriting its data with the 2022-Q4 numbers.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The template sheet had 3 fund rows (rows 2-4); 2022-Q4 only has 1, so
# drop rows 3 and 4 entirely (also shrinks the sheet dimension to A1:H2).
$newSheet.Range("A3:H4").EntireRow.Delete()

# Header row: "基金金额" -> "基金规模" (the rest of the header is unchanged).
$newSheet.Range("D1").Value = "基金规模"

# Fill in the single 2022-Q4 fund row. Codes/ratios are text in this
# workbook (matching the sibling quarter sheets), so force text with a
# leading apostrophe and then drop the resulting "quote prefix" style so
# the cell ends up with no explicit style, same as its neighbours.
$newSheet.Range("B2").Value = "'001068"
$newSheet.Range("B2").ClearFormats()
$newSheet.Range("C2").Value = "国新国证新锐灵活配置混合"
$newSheet.Range("D2").Value = "'0.21"
$newSheet.Range("D2").ClearFormats()
$newSheet.Range("E2").Value = "'75.37"
$newSheet.Range("E2").ClearFormats()
$newSheet.Range("F2").Value = "'3.13"
$newSheet.Range("F2").ClearFormats()
$newSheet.Range("G2").Value = "'0.0066"
$newSheet.Range("G2").ClearFormats()
$newSheet.Range("H2").Value = 8

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row right under the
#    header for the 2022-Q4 quarter, pushing the existing quarters down
#    one row (their data is unchanged, only the running index in column A
#    is renumbered to stay sequential).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Range("A2").EntireRow.Insert()

# Give the new row the same styling as the other index/data rows.
$summary.Range("A5").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B3:D3").Copy()
$summary.Range("B2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.01

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ---------------------------------------------------------------------
# 3) Restore the original active sheet ("2020-Q4", the last tab).
# ---------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
